$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells keep their text formatting (e.g. "1.00", "0.590")
# instead of being auto-converted to numbers when values are assigned.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '64.746.84'
$ws.Range("E2").Value = '  -0.29%  '

# Row 3
$ws.Range("D3").Value = '3.164.45'
$ws.Range("E3").Value = '  -2.28%  '

# Row 4
$ws.Range("E4").Value = '  +0.06%  '

# Row 5
$ws.Range("D5").Value = '570.71'
$ws.Range("E5").Value = '  -1.11%  '

# Row 6
$ws.Range("D6").Value = '165.21'
$ws.Range("E6").Value = '  -5.43%  '

# Row 7
$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D7").Value = '0.590'
$ws.Range("E7").Value = '  -6.21%  '

# Row 8
$ws.Range("B8").Value = 'USDC'
$ws.Range("C8").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.07%  '

# Row 9
$ws.Range("D9").Value = '0.118'
$ws.Range("E9").Value = '  -3.74%  '

# Row 10
$ws.Range("D10").Value = '6.69'
$ws.Range("E10").Value = '  -1.45%  '

# Row 11
$ws.Range("D11").Value = '0.385'
$ws.Range("E11").Value = '  -1.46%  '

# Row 12
$ws.Range("D12").Value = '3.723.24'

# Row 14
$ws.Range("D14").Value = '64.542.68'
$ws.Range("E14").Value = '  -0.78%  '

# Row 15
$ws.Range("D15").Value = '25.37'
$ws.Range("E15").Value = '  -1.02%  '

# Row 16
$ws.Range("D16").Value = '3.168.49'
$ws.Range("E16").Value = '  -2.08%  '

# Row 17
$ws.Range("D17").Value = '0.0000156'
$ws.Range("E17").Value = '  -2.02%  '

# Row 18
$ws.Range("D18").Value = '414.39'
$ws.Range("E18").Value = '  +0.25%  '

# Row 19
$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D19").Value = '12.65'

# Row 20
$ws.Range("B20").Value = 'Polkadot'
$ws.Range("C20").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D20").Value = '5.28'
$ws.Range("E20").Value = '  -1.98%  '

# Row 21
$ws.Range("D21").Value = '7.14'
$ws.Range("E21").Value = '  -0.90%  '

# Row 22
$ws.Range("E22").Value = '  +0.09%  '

# Row 23
$ws.Range("D23").Value = '68.15'
$ws.Range("E23").Value = '  -3.42%  '

# Row 24
$ws.Range("E24").Value = '  -1.62%  '

# Row 25
$ws.Range("E25").Value = '  -1.95%  '

# Row 26
$ws.Range("E26").Value = '  -6.44%  '

# Row 27
$ws.Range("D27").Value = '8.91'
$ws.Range("E27").Value = '  -2.39%  '

# Row 28
$ws.Range("D28").Value = '0.998'
$ws.Range("E28").Value = '  -0.19%  '

# Row 29
$ws.Range("D29").Value = '1.82'
$ws.Range("E29").Value = '  -2.94%  '

# Row 30
$ws.Range("D30").Value = '21.28'
$ws.Range("E30").Value = '  -2.51%  '

# Row 31
$ws.Range("B31").Value = 'NEARProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D31").Value = '4.93'
$ws.Range("E31").Value = '  -1.60%  '

# Row 32
$ws.Range("B32").Value = 'Aptos'
$ws.Range("C32").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D32").Value = '6.33'
$ws.Range("E32").Value = '  -1.60%  '

# Row 33
$ws.Range("E33").Value = '  -2.10%  '

# Row 34
$ws.Range("D34").Value = '155.57'
$ws.Range("E34").Value = '  -0.36%  '

# Row 35
$ws.Range("D35").Value = '1.35'
$ws.Range("E35").Value = '  -3.30%  '

# Row 36
$ws.Range("D36").Value = '2.722.02'
$ws.Range("E36").Value = '  -3.81%  '

# Row 37
$ws.Range("E37").Value = '  -2.92%  '

# Row 38
$ws.Range("B38").Value = 'EnergySwap'
$ws.Range("C38").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D38").Value = '23.64'
$ws.Range("E38").Value = '  -7.19%  '

# Row 39
$ws.Range("B39").Value = 'Filecoin'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D39").Value = '4.08'
$ws.Range("E39").Value = '  -3.12%  '

# Row 40
$ws.Range("D40").Value = '0.702'
$ws.Range("E40").Value = '  -3.91%  '

# Row 41
$ws.Range("D41").Value = '0.0634'
$ws.Range("E41").Value = '  +1.15%  '

# Row 42
$ws.Range("D42").Value = '5.54'
$ws.Range("E42").Value = '  -3.74%  '

# Row 43
$ws.Range("D43").Value = '0.0262'
$ws.Range("E43").Value = '  -0.46%  '

# Row 44
$ws.Range("D44").Value = '289.08'
$ws.Range("E44").Value = '  -5.64%  '

# Row 45
$ws.Range("D45").Value = '21.24'
$ws.Range("E45").Value = '  -4.16%  '

# Row 46
$ws.Range("D46").Value = '1.00'
$ws.Range("E46").Value = '  +0.09%  '

# Row 47
$ws.Range("E47").Value = '  -2.95%  '

# Row 48
$ws.Range("D48").Value = '1.95'
$ws.Range("E48").Value = '  -11.87%  '

# Row 49
$ws.Range("D49").Value = '10.45'
$ws.Range("E49").Value = '  +0.59%  '

# Row 50
$ws.Range("D50").Value = '5.76'
$ws.Range("E50").Value = '  -0.87%  '

# Row 51
$ws.Range("D51").Value = '0.901'
$ws.Range("E51").Value = '  -3.89%  '

# Restore default (General) style on the Price column so unchanged cells
# keep their original appearance/style index.
$ws.Range("D2:D51").Style = "Normal"
